$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 13's data)
$ws.Range("D2").Value = 44553
$ws.Range("H2").Value = 'Inferno'
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 45000
$ws.Range("L2").Value = 45000
$ws.Range("M2").Value = 45000
$ws.Range("P2").Value = 1800

# Row 3 (was row 6's data)
$ws.Range("D3").Value = 44425
$ws.Range("H3").Value = 'Americana (o)'
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 75000
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = 75000
$ws.Range("N3").Value = '$/caja 25 kilos'
$ws.Range("P3").Value = 3000
$ws.Range("Q3").Value = 25

# Row 4 (was row 8's data)
$ws.Range("D4").Value = 44340
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 35000
$ws.Range("L4").Value = 35000
$ws.Range("M4").Value = 35000
$ws.Range("P4").Value = 1400

# Row 5 (was row 17's data)
$ws.Range("D5").Value = 44449
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 80000
$ws.Range("L5").Value = 80000
$ws.Range("M5").Value = 80000
$ws.Range("N5").Value = '$/caja 25 kilos'
$ws.Range("P5").Value = 3200
$ws.Range("Q5").Value = 25

# Row 6 (was row 18's data)
$ws.Range("D6").Value = 44449
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 20
$ws.Range("N6").Value = '$/caja 15 kilos'
$ws.Range("P6").Value = 5000
$ws.Range("Q6").Value = 15

# Row 7 (was row 16's data)
$ws.Range("D7").Value = 44544
$ws.Range("H7").Value = 'Inferno'
$ws.Range("J7").Value = 12
$ws.Range("K7").Value = 35000
$ws.Range("L7").Value = 35000
$ws.Range("M7").Value = 35000
$ws.Range("P7").Value = 1400

# Row 8 (was row 7's data)
$ws.Range("D8").Value = 44319
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 30000
$ws.Range("P8").Value = 1200

# Row 9 (was row 15's data)
$ws.Range("D9").Value = 44474
$ws.Range("J9").Value = 18
$ws.Range("K9").Value = 100000
$ws.Range("L9").Value = 100000
$ws.Range("M9").Value = 100000
$ws.Range("P9").Value = 4000

# Row 10 (was row 4's data)
$ws.Range("D10").Value = 44343
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 36000
$ws.Range("L10").Value = 36000
$ws.Range("M10").Value = 36000
$ws.Range("P10").Value = 1440

# Row 11 (was row 9's data)
$ws.Range("D11").Value = 44460
$ws.Range("I11").Value = 'Primera'
$ws.Range("K11").Value = 95000
$ws.Range("L11").Value = 95000
$ws.Range("M11").Value = 95000
$ws.Range("P11").Value = 3800

# Row 12 (was row 10's data)
$ws.Range("D12").Value = 44326
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("P12").Value = 1200

# Row 13 (was row 11's data)
$ws.Range("D13").Value = 44581
$ws.Range("H13").Value = 'Americana (o)'
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 17000
$ws.Range("P13").Value = 680

# Row 14 (was row 2's data)
$ws.Range("D14").Value = 44446
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 78000
$ws.Range("L14").Value = 78000
$ws.Range("M14").Value = 78000
$ws.Range("P14").Value = 3120

# Row 15 (was row 3's data)
$ws.Range("D15").Value = 44446
$ws.Range("H15").Value = 'Inferno'
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 80000
$ws.Range("L15").Value = 80000
$ws.Range("M15").Value = 80000
$ws.Range("N15").Value = '$/caja 15 kilos'
$ws.Range("P15").Value = 5333
$ws.Range("Q15").Value = 15

# Row 16 (was row 14's data)
$ws.Range("D16").Value = 44421
$ws.Range("H16").Value = 'Americana (o)'
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 75000
$ws.Range("L16").Value = 75000
$ws.Range("M16").Value = 75000
$ws.Range("P16").Value = 3000

# Row 17 (was row 5's data)
$ws.Range("D17").Value = 44193
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 46000
$ws.Range("L17").Value = 46000
$ws.Range("M17").Value = 46000
$ws.Range("N17").Value = '$/caja 15 kilos'
$ws.Range("P17").Value = 3067
$ws.Range("Q17").Value = 15

# Row 18 (was row 12's data)
$ws.Range("D18").Value = 44221
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24545
$ws.Range("N18").Value = '$/caja 25 kilos'
$ws.Range("P18").Value = 982
$ws.Range("Q18").Value = 25
